# Apply "merged input from Wanting" edit to slide 3 (Content Placeholder 2)
#
# 1) "Symptom is now a list instead of a container." ->
#    "Symptom " / "and annotator " / "is now a list instead of a container."
#    (split into three runs)
#
# 2) "-> Work on " / "example implementation " / "in IETF 121 hackathon."
#    -> merged back into a single run
#    "-> Work on example implementation in IETF 121 hackathon."

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Part 1: split the "Symptom ..." sentence into three runs ---------

$symptomRun = $tr.Find("Symptom ")
$symptomSub = $tr.Characters($symptomRun.Start, $symptomRun.Length)
$symptomSub.Text = "Symptom "

$tailRun = $tr.Find("is now a list instead of a container.")
$tailSub = $tr.Characters($tailRun.Start, $tailRun.Length)
$tailSub.Text = "and annotator is now a list instead of a container."

$tailRun2 = $tr.Find("is now a list instead of a container.")
$tailSub2 = $tr.Characters($tailRun2.Start, $tailRun2.Length)
$tailSub2.Text = "is now a list instead of a container."

# --- Part 2: merge the "Work on / example implementation / in IETF ..." runs --

$workOn = $tr.Find("-> Work on example implementation in IETF 121 hackathon.")
$workOn.Text = "-> Work on example implementation in IETF 121 hackathon."
